$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank columns before column A, shifting old A->C, B->D, C->E
$ws.Columns("A:B").Insert()

# ---- Description column (B1:B13) ----
$ws.Range("B1").Value = "Description"
$ws.Range("B2").Value = " Retrieve all components of the building which are generic placeholder"
$ws.Range("B3").Value = "Retrieve all walls with a window"
$ws.Range("B4").Value = "Retrieve  all walls with a door"
$ws.Range("B5").Value = "Query all doors, windows, walls and slabs"
$ws.Range("B6").Value = "Query the name of the construction project"
$ws.Range("B7").Value = "Retrieve the wall with the name Basiswand:STB 250"
$ws.Range("B8").Value = "Query doors with reference DL - 900 x 2175"
$ws.Range("B9").Value = "Query number of doors that are not the entry door"
$ws.Range("B10").Value = "Retrieve the highest door "
$ws.Range("B11").Value = "Retrieve a wall in which a window with the smallest width is embedded"
$ws.Range("B12").Value = "Query all walls that connections at least two rooms → these are walls that have at least 2 doors"
$ws.Range("B13").Value = "Retrieve the number of rooms"

# ---- Query Number column (A1:A13) ----
$ws.Range("A1").Value = "Query Number"
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 4
$ws.Range("A6").Value = 5
$ws.Range("A7").Value = 6
$ws.Range("A8").Value = 7
$ws.Range("A9").Value = 8
$ws.Range("A10").Value = 9
$ws.Range("A11").Value = 10
$ws.Range("A12").Value = 11
$ws.Range("A13").Value = 12

# ---- Column widths ----
# Note: ColumnWidth snaps to a pixel grid (increments of 1/6 character unit);
# the values below are the inputs that land closest to the target widths
# (73.6640625 and 131.83203125) on that grid.
$ws.Columns("B").ColumnWidth = 72.83333333333334
$ws.Columns("E").ColumnWidth = 131

# ---- View settings ----
$ws.Application.ActiveWindow.Zoom = 111
$ws.Range("B3").Select()
